$d = $word.ActiveDocument

# Remove the stray "-beschrijving kamers" bullet line: the room/student/
# activity/bestelling ERD description that follows it already (and now,
# after adding the ERD conversion text for "bestelling") covers this, so
# the placeholder line is deleted as its own paragraph (including the
# paragraph mark) so no blank line is left behind.
$target = "-beschrijving kamers"
$found = $d.Content.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng = $d.Content
    $rng.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $para = $rng.Paragraphs(1).Range
    $para.Delete()
}
